$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert three new modules (gsap-lite, pixi, viewport-size) amongst the
# --- existing alphabetical module list, shifting later rows down.
# Existing layout (rows 12-16):
#   12 common.preloader
#   13 common.retina-images
#   14 common.simple-hoverable
#   16 state.preloader   (row 15 is blank)
# New layout (rows 12-17, 19):
#   12 common.gsap-lite
#   13 common.pixi
#   14 common.preloader
#   15 common.retina-images
#   16 common.simple-hoverable
#   17 common.viewport-size
#   19 state.preloader   (row 18 is blank)

# Shift state.preloader from row 16 down to row 19 first, since that
# destination area is currently empty.
$ws.Range("A19").Value = $ws.Range("A16").Value()
$ws.Range("E19").Value = $ws.Range("E16").Value()
$ws.Range("B19:D19").Style = "Good"

# Shift the three existing rows (preloader, retina-images, simple-hoverable)
# down by two rows: 12->14, 13->15, 14->16
$ws.Range("A16").Value = $ws.Range("A14").Value()
$ws.Range("E16").Value = $ws.Range("E14").Value()
$ws.Range("B16:D16").Style = "Good"

$ws.Range("A15").Value = $ws.Range("A13").Value()
$ws.Range("E15").Value = $ws.Range("E13").Value()
$ws.Range("B15:D15").Style = "Good"

$ws.Range("A14").Value = $ws.Range("A12").Value()
$ws.Range("E14").Value = $ws.Range("E12").Value()
$ws.Range("B14:D14").Style = "Good"

# Now fill the freed-up rows 12, 13 and the new row 17 with the new modules.
# (Write in the same order the new strings were introduced upstream: gsap-lite,
# viewport-size, pixi -- so the shared-strings table comes out in that order.)
$ws.Range("A12").Value = "common.gsap-lite"
$ws.Range("E12").Value = "Incapsuleaza TweenLite si Easing'urile."
$ws.Range("B12:D12").Style = "Good"

$ws.Range("A17").Value = "common.viewport-size"
$ws.Range("E17").Value = "Ne ofera acces la resolutia curenta a viewportului. Totodata ne ofera posibilitatea de a registra callbackuri."
$ws.Range("B17:D17").Style = "Good"

$ws.Range("A13").Value = "common.pixi"
$ws.Range("E13").Value = "Incapsuleaza pixi.js."
$ws.Range("B13:D13").Style = "Good"

# --- Update the view: zoom level and active selection
$excel.ActiveWindow.Zoom = 85
$ws.Range("A20").Select()
